$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the three new rows of notes at B10:B12 (row 9 left empty, matching original).
# Assign values in this order so that new entries land in the shared-strings
# table in the same sequence as the target workbook.
$ws.Range("B11").Value = "31,32,34 utilizzati per convivenze di fatto"
$ws.Range("B12").Value = "98 utilizzato in entrambi i casi"
$ws.Range("B10").Value = "1,2,3 utilizzati per unioni civili"

# Update the selected cell/window view to match new active cell
$ws.Range("B16").Select()
